$d = $word.ActiveDocument

$d.Content.Find.Execute("   * RSP ITEM", $false, $false, $false, $false, $false, $true, 1, $false, "   * HQ Power - RSP", 2)
$d.Content.Find.Execute("   * sdgaty,", $false, $false, $false, $false, $false, $true, 1, $false, "   * TR/2025/12", 2)
$d.Content.Find.Execute("   * 28-07-2025", $false, $false, $false, $false, $false, $true, 1, $false, "   * 30-07-2025", 2)
$d.Content.Find.Execute("   * 07-28-2025", $false, $false, $false, $false, $false, $true, 1, $false, "   * 07-30-2025", 2)
